# Scheduled runner update: refresh cached market-board pricing / profit
# figures across the Leve profit sheets (H:N columns) per latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 637.4054
$ws.Range("I12").Value = 402.5
$ws.Range("J12").Value = 3299.6667
$ws.Range("K12").Value = 402.5
$ws.Range("L12").Value = 3299.6667
$ws.Range("M12").Value = -232.5
$ws.Range("N12").Value = -3639.6667
$ws.Range("H19").Value = 1744.875
$ws.Range("J19").Value = 2247.25
$ws.Range("L19").Value = 2247.25
$ws.Range("N19").Value = -2597.25
$ws.Range("H107").Value = 1486.0741
$ws.Range("I107").Value = 1058.1765
$ws.Range("K107").Value = 1058.1765
$ws.Range("M107").Value = 861.8235
$ws.Range("H137").Value = 2693.9033
$ws.Range("I137").Value = 920.9375
$ws.Range("J137").Value = 4585.067
$ws.Range("K137").Value = 2762.8125
$ws.Range("L137").Value = 13755.201
$ws.Range("M137").Value = -212.8125
$ws.Range("N137").Value = -18855.201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 6882.909
$ws.Range("I110").Value = 4541.6
$ws.Range("J110").Value = 7571.5293
$ws.Range("K110").Value = 4541.6
$ws.Range("L110").Value = 7571.5293
$ws.Range("M110").Value = -2496.6
$ws.Range("N110").Value = -11661.5293
$ws.Range("H132").Value = 561840.4399999999
$ws.Range("I132").Value = 723963.6
$ws.Range("J132").Value = 51152.4
$ws.Range("K132").Value = 2171890.8
$ws.Range("L132").Value = 153457.2
$ws.Range("M132").Value = -2169360.8
$ws.Range("N132").Value = -158517.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3333529
$ws.Range("I22").Value = 293.5
$ws.Range("J22").Value = 10000000
$ws.Range("K22").Value = 293.5
$ws.Range("L22").Value = 10000000
$ws.Range("M22").Value = -120.5
$ws.Range("N22").Value = -10000346
$ws.Range("H94").Value = 4371.8823
$ws.Range("I94").Value = 3293.24
$ws.Range("J94").Value = 7368.1113
$ws.Range("K94").Value = 3293.24
$ws.Range("L94").Value = 7368.1113
$ws.Range("M94").Value = -2842.24
$ws.Range("N94").Value = -8270.1113
$ws.Range("H96").Value = 3607.5
$ws.Range("I96").Value = 3607.5
$ws.Range("K96").Value = 3607.5
$ws.Range("M96").Value = -861.5
$ws.Range("H99").Value = 8987.513999999999
$ws.Range("I99").Value = 9166.297
$ws.Range("J99").Value = 8798.514999999999
$ws.Range("K99").Value = 9166.297
$ws.Range("L99").Value = 8798.514999999999
$ws.Range("M99").Value = -7668.297
$ws.Range("N99").Value = -11794.515
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8847.254999999999
$ws.Range("I31").Value = 12123.056
$ws.Range("J31").Value = 7253.6216
$ws.Range("K31").Value = 12123.056
$ws.Range("L31").Value = 7253.6216
$ws.Range("M31").Value = -11828.056
$ws.Range("N31").Value = -7843.6216
$ws.Range("H34").Value = 8847.254999999999
$ws.Range("I34").Value = 12123.056
$ws.Range("J34").Value = 7253.6216
$ws.Range("K34").Value = 12123.056
$ws.Range("L34").Value = 7253.6216
$ws.Range("M34").Value = -11921.056
$ws.Range("N34").Value = -7657.6216
$ws.Range("H58").Value = 7416.9644
$ws.Range("I58").Value = 4594.952
$ws.Range("J58").Value = 15883
$ws.Range("K58").Value = 4594.952
$ws.Range("L58").Value = 15883
$ws.Range("M58").Value = -4391.952
$ws.Range("N58").Value = -16289
$ws.Range("H99").Value = 4633086.5
$ws.Range("J99").Value = 3825.7
$ws.Range("L99").Value = 3825.7
$ws.Range("N99").Value = -6821.7
$ws.Range("H122").Value = 3672.5483
$ws.Range("I122").Value = 2972.05
$ws.Range("K122").Value = 8916.150000000001
$ws.Range("M122").Value = -6466.150000000001
$ws.Range("H126").Value = 4633086.5
$ws.Range("J126").Value = 3825.7
$ws.Range("L126").Value = 11477.1
$ws.Range("N126").Value = -16417.1
$ws.Range("H132").Value = 6269.2036
$ws.Range("I132").Value = 4736.606
$ws.Range("J132").Value = 8677.571
$ws.Range("K132").Value = 14209.818
$ws.Range("L132").Value = 26032.713
$ws.Range("M132").Value = -11679.818
$ws.Range("N132").Value = -31092.713
$ws.Range("H136").Value = 7416.9644
$ws.Range("I136").Value = 4594.952
$ws.Range("J136").Value = 15883
$ws.Range("K136").Value = 13784.856
$ws.Range("L136").Value = 47649
$ws.Range("M136").Value = -11234.856
$ws.Range("N136").Value = -52749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5245.104
$ws.Range("I102").Value = 4224.643
$ws.Range("J102").Value = 6673.75
$ws.Range("K102").Value = 4224.643
$ws.Range("L102").Value = 6673.75
$ws.Range("M102").Value = -2602.643
$ws.Range("N102").Value = -9917.75
$ws.Range("H113").Value = 8812.182000000001
$ws.Range("I113").Value = 7330
$ws.Range("J113").Value = 9368
$ws.Range("K113").Value = 7330
$ws.Range("L113").Value = 9368
$ws.Range("M113").Value = -5160
$ws.Range("N113").Value = -13708
$ws.Range("H126").Value = 25007954
$ws.Range("I126").Value = 41667720
$ws.Range("J126").Value = 18301.125
$ws.Range("K126").Value = 125003160
$ws.Range("L126").Value = 54903.375
$ws.Range("M126").Value = -125000690
$ws.Range("N126").Value = -59843.375
$ws.Range("H132").Value = 3757.7073
$ws.Range("I132").Value = 3658.889
$ws.Range("K132").Value = 10976.667
$ws.Range("M132").Value = -8446.667000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8465.763000000001
$ws.Range("I7").Value = 7289.346
$ws.Range("J7").Value = 11014.667
$ws.Range("K7").Value = 7289.346
$ws.Range("L7").Value = 11014.667
$ws.Range("M7").Value = -7177.346
$ws.Range("N7").Value = -11238.667
$ws.Range("H40").Value = 8994.666999999999
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H100").Value = 4472.385
$ws.Range("I100").Value = 4831
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 4831
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -4290
$ws.Range("N100").Value = -3582
$ws.Range("H122").Value = 4998951
$ws.Range("I122").Value = 9984898
$ws.Range("J122").Value = 13004
$ws.Range("K122").Value = 29954694
$ws.Range("L122").Value = 39012
$ws.Range("M122").Value = -29952244
$ws.Range("N122").Value = -43912
$ws.Range("H126").Value = 8465.763000000001
$ws.Range("I126").Value = 7289.346
$ws.Range("J126").Value = 11014.667
$ws.Range("K126").Value = 21868.038
$ws.Range("L126").Value = 33044.001
$ws.Range("M126").Value = -19398.038
$ws.Range("N126").Value = -37984.001
$ws.Range("H132").Value = 8065.8945
$ws.Range("I132").Value = 7257.75
$ws.Range("K132").Value = 21773.25
$ws.Range("M132").Value = -19243.25
$ws.Range("H136").Value = 45465340
$ws.Range("I136").Value = 6905.6924
$ws.Range("K136").Value = 20717.0772
$ws.Range("M136").Value = -18167.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1694.875
$ws.Range("I100").Value = 1704.75
$ws.Range("J100").Value = 1675.125
$ws.Range("K100").Value = 3409.5
$ws.Range("L100").Value = 3350.25
$ws.Range("M100").Value = -2868.5
$ws.Range("N100").Value = -4432.25
$ws.Range("H107").Value = 16667855
$ws.Range("I107").Value = 20001206
$ws.Range("K107").Value = 60003618
$ws.Range("M107").Value = -60001698
$ws.Range("H122").Value = 8144.3096
$ws.Range("I122").Value = 3263.3872
$ws.Range("J122").Value = 21899.637
$ws.Range("K122").Value = 9790.161599999999
$ws.Range("L122").Value = 65698.91099999999
$ws.Range("M122").Value = -7340.161599999999
$ws.Range("N122").Value = -70598.91099999999
$ws.Range("H126").Value = 3051.875
$ws.Range("I126").Value = 1811.4814
$ws.Range("K126").Value = 5434.4442
$ws.Range("M126").Value = -2964.4442
$ws.Range("H132").Value = 6413.056
$ws.Range("I132").Value = 4935.5156
$ws.Range("K132").Value = 14806.5468
$ws.Range("M132").Value = -12276.5468
